# Update the NBA Top10 leader sheets (leaders_nba.xlsx)
# The "Valor" column stores numeric-looking values as TEXT (shared strings),
# so we write them with a leading apostrophe to force text, then reset the
# cell style back to Normal so no stray number-format style is introduced.
$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- Sheet: Arremessos de 3 Pontos (3-point makes) ---
$ws = $wb.Worksheets.Item("Arremessos de 3 Pontos")
Set-TextValue $ws.Range("D3") "3.9"
$ws.Range("A6").Value = 5
Set-TextValue $ws.Range("D6") "3.3"

# --- Sheet: Assistências (Assists) ---
$ws = $wb.Worksheets.Item("Assistências")
Set-TextValue $ws.Range("D3") "10.9"
Set-TextValue $ws.Range("D4") "9.5"
Set-TextValue $ws.Range("D5") "8.9"
Set-TextValue $ws.Range("D6") "8.4"

# --- Sheet: Pontos (Points) ---
$ws = $wb.Worksheets.Item("Pontos")
Set-TextValue $ws.Range("D3") "34.2"
Set-TextValue $ws.Range("D5") "30.8"
$ws.Range("B6").Value = "Donovan Mitchell"
$ws.Range("C6").Value = "CLE"
Set-TextValue $ws.Range("D6") "28.4"

# --- Sheet: Rebotes (Rebounds) ---
$ws = $wb.Worksheets.Item("Rebotes")
$ws.Range("B4").Value = "Anthony Davis"
$ws.Range("C4").Value = "LAL"
$ws.Range("A5").Value = 4
Set-TextValue $ws.Range("D5") "12.0"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Jalen Duren"
$ws.Range("C6").Value = "DET"
Set-TextValue $ws.Range("D6") "12.0"

# --- Sheet: Roubos (Steals) ---
$ws = $wb.Worksheets.Item("Roubos")
$ws.Range("A4").Value = 2
Set-TextValue $ws.Range("D4") "1.9"

# --- Sheet: Tocos (Blocks) ---
$ws = $wb.Worksheets.Item("Tocos")
Set-TextValue $ws.Range("D2") "3.2"
Set-TextValue $ws.Range("D6") "2.5"
